# "empty option in mapping fields"
#
# 1) Color the "in step 1, when file imported successfully, ..." bullet
#    the same green (00A933 / wdColor 3385600) used by the other TO-DO
#    items - applies to both the run text and the paragraph mark.
# 2) Normal style: stop allowing punctuation to overflow the text
#    extents (w:overflowPunct false), exposed via
#    ParagraphFormat.HangingPunctuation on the style.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*in step 1, when file imported successfully*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Font.Color = 3385600
}

$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = $false
